$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-12 Tuesday" "2023-09-13 Wednesday"

Replace-Text "20×27=" "60×65="
Replace-Text "48×26=" "95×81="
Replace-Text "89×98=" "50×29="
Replace-Text "99×91=" "98×78="
Replace-Text "95×91=" "28×38="

Replace-Text "42×74=" "46×33="
Replace-Text "76×63=" "32×13="
Replace-Text "82×90=" "64×60="
Replace-Text "95×82=" "90×43="
Replace-Text "36×91=" "28×58="

Replace-Text "94×51=" "32×41="
Replace-Text "78×18=" "30×15="
Replace-Text "97×38=" "66×62="
Replace-Text "54×21=" "12×92="
Replace-Text "71×69=" "76×96="

Replace-Text "53×11=" "60×59="
Replace-Text "43×98=" "32×93="
Replace-Text "87×93=" "97×38="
Replace-Text "89×78=" "83×44="
Replace-Text "74×95=" "73×23="

Replace-Text "49×54=" "49×92="
Replace-Text "55×96=" "35×81="
Replace-Text "96×34=" "99×71="
Replace-Text "84×49=" "20×72="
Replace-Text "20×90=" "22×47="
